# Changed the localize handlebars function to expose the data model and the
# calculates: every "{{member_name}}" / "{{household_id}}" handlebars
# reference on the "survey" sheet becomes "{{data.member_name}}" /
# "{{data.household_id}}" so the template also has access to the full data
# model (matches the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Range("D2").Value = "Data for household: {{data.household_id}}"
$ws.Range("D4").Value = "Enter age of {{data.member_name}}:"
$ws.Range("D5").Value = "Enter sex of {{data.member_name}}:"
$ws.Range("D6").Value = "{{data.member_name}} age is {{evaluate calculates.ageIsOddOrEven}} in {{setting 'table_id'}}"
$ws.Range("D7").Value = "Does {{data.member_name}} contribute to the household income?"

$ws.Activate()
$ws.Range("D6").Select()
